# Deploying to gh-pages: add the 2023 column (T) to the suicide-mortality table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting of column S into the new column T (rows 4-14) ---
$ws.Range("S4:S14").Copy($ws.Range("T4:T14"))

# --- 2. Write the new 2023 figures into column T ---
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 4.8
$ws.Range("T6").Value = 5.7
$ws.Range("T7").Value = 1.9
$ws.Range("T8").Value = 8.9
$ws.Range("T9").Value = 11.9
$ws.Range("T10").Value = 2.5
$ws.Range("T11").Value = 0.7
$ws.Range("T12").Value = 12.7
$ws.Range("T13").Value = 1.1
$ws.Range("T14").Value = 2.2

# --- 3. Resize columns A:C and D to their new widths ---
$ws.Range("A1:C1").ColumnWidth = 30.666666666666668
$ws.Range("D1").ColumnWidth = 8.833333333333334

# --- 4. Adjust row heights for the header rows + the first data row ---
$ws.Range("A1").RowHeight = 30
$ws.Range("A2").RowHeight = 15.75
$ws.Range("A5").RowHeight = 16.5

# --- 5. Reset the selection back to the default (A1) ---
$ws.Range("A1").Select() | Out-Null

Write-Output "edit applied"
